# spade_Zambia_wom_vita.xlsx - "updated output" edit
# - refresh the run Start_time / End_time stamps on the "Info" sheet
# - bump a few R package version numbers on the "sessionInfo" sheet
# - drop the "backports" row from the Loaded_only package table
#   (the row below it, "boot", shifts up to take its place and the
#   now-unused last row of that table is cleared)

$wb = $excel.ActiveWorkbook
$wsInfo    = $wb.Worksheets.Item("Info")
$wsSession = $wb.Worksheets.Item("sessionInfo")

# --- Info sheet: Start_time (A26/B26) / End_time (A27/B27) ---
$wsInfo.Range("B26").Value = "Thu Nov 19 15:49:32 2020"
$wsInfo.Range("B27").Value = "Thu Nov 19 15:49:44 2020"

# --- sessionInfo sheet: package version bumps ---
$wsSession.Range("G2").Value  = "1.0.0"   # here: 0.1 -> 1.0.0
$wsSession.Range("J3").Value  = "2.0.1"   # magrittr: 1.5 -> 2.0.1
$wsSession.Range("J10").Value = "2.0.2"   # rprojroot: 1.3-2 -> 2.0.2

# --- sessionInfo sheet: remove the "backports" row from the
#     Loaded_only table (I/J columns). "boot"/"1.3-25" (row 16)
#     moves up into row 15, and row 16's I/J cells are cleared.
$wsSession.Range("I15").Value = "boot"
$wsSession.Range("J15").Value = "1.3-25"
$wsSession.Range("I16").ClearContents()
$wsSession.Range("J16").ClearContents()
